$wb = $excel.ActiveWorkbook

# The "Italy" sheet is the 7th sheet and is currently the active/selected tab.
$wsItaly = $wb.Worksheets.Item(7)

# Create the new "Spain" sheet as a copy of "Italy", placed immediately after it.
$wsItaly.Copy($null, $wsItaly)
$wsSpain = $wb.Worksheets.Item(8)
$wsSpain.Name = "Spain"

# Update the market-specific values on the new sheet.
$wsSpain.Range("B2").Value = "Spain Market"
$wsSpain.Range("B4").Value = "NGC-3103/T2035/"

# Resize columns B:D to match the final layout of the Spain sheet.
$wsSpain.Columns.Item(2).ColumnWidth = 21.67
$wsSpain.Columns.Item(3).ColumnWidth = 16.11
$wsSpain.Columns.Item(4).ColumnWidth = 33.56

# Rows 3-5 grow to a two-line height on the new sheet.
$wsSpain.Rows.Item(3).RowHeight = 28.8
$wsSpain.Rows.Item(4).RowHeight = 28.8
$wsSpain.Rows.Item(5).RowHeight = 28.8

# The Italy sheet is no longer the active tab; its selection becomes the full
# used range instead of the previous B4 cell selection.
$wsItaly.Range("A1:D21").Select()

# Make the new Spain sheet the active tab with its original selection state.
$wsSpain.Activate()
$wsSpain.Range("B4").Select()
